$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.193.00'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '3.913.88'
$ws.Range("E3").Value = '  +2.84%  '
$ws.Range("E4").Value = '  +0.28%  '
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '468.11'
$c.Style = $origStyle
$ws.Range("E5").Value = '  +7.59%  '
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '144.83'
$c.Style = $origStyle
$ws.Range("E6").Value = '  +4.18%  '
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.621'
$c.Style = $origStyle
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("E8").Value = '  -0.13%  '
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.730'
$c.Style = $origStyle
$ws.Range("E9").Value = '  -1.33%  '
$ws.Range("E10").Value = '  +7.07%  '
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0000341'
$c.Style = $origStyle
$ws.Range("E11").Value = '  +5.22%  '
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '43.13'
$c.Style = $origStyle
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").Value = '4.538.12'
$ws.Range("E13").Value = '  +3.11%  '
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '10.33'
$c.Style = $origStyle
$ws.Range("E14").Value = '  -1.64%  '
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.90'
$c.Style = $origStyle
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '3.919.76'
$ws.Range("E16").Value = '  +2.11%  '
$ws.Range("E17").Value = '  -0.38%  '
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '19.80'
$c.Style = $origStyle
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = '67.381.65'
$ws.Range("E20").Value = '  +1.07%  '
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '430.08'
$c.Style = $origStyle
$ws.Range("E21").Value = '  +4.12%  '
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.36'
$c.Style = $origStyle
$ws.Range("E22").Value = '  +3.37%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.56'
$c.Style = $origStyle
$ws.Range("E23").Value = '  -1.75%  '
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '87.45'
$c.Style = $origStyle
$ws.Range("E24").Value = '  +2.52%  '
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.54'
$c.Style = $origStyle
$ws.Range("E25").Value = '  +5.16%  '
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '38.39'
$c.Style = $origStyle
$ws.Range("E26").Value = '  +3.36%  '
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '10.23'
$c.Style = $origStyle
$ws.Range("E27").Value = '  +3.47%  '
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.73'
$c.Style = $origStyle
$ws.Range("E28").Value = '  +2.64%  '
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.58'
$c.Style = $origStyle
$ws.Range("E29").Value = '  -1.34%  '
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '727.82'
$c.Style = $origStyle
$ws.Range("E30").Value = '  +1.90%  '
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '13.51'
$c.Style = $origStyle
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("E32").Value = '  -3.27%  '
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.82'
$c.Style = $origStyle
$ws.Range("E33").Value = '  +2.06%  '
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '42.82'
$c.Style = $origStyle
$ws.Range("E34").Value = '  +2.66%  '
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.154'
$c.Style = $origStyle
$ws.Range("E35").Value = '  +2.68%  '
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '57.85'
$c.Style = $origStyle
$ws.Range("E36").Value = '  +3.08%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0799'
$ws.Range("E37").Value = '  +13.76%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = $origStyle
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E39").Value = '  -4.22%  '
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0477'
$c.Style = $origStyle
$ws.Range("E40").Value = '  +0.41%  '
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.04'
$c.Style = $origStyle
$ws.Range("E41").Value = '  +4.71%  '
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.60'
$c.Style = $origStyle
$ws.Range("E42").Value = '  -5.60%  '
$ws.Range("E43").Value = '  -1.13%  '
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.336'
$c.Style = $origStyle
$ws.Range("E44").Value = '  +3.18%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("E46").Value = '  +3.79%  '
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.16'
$c.Style = $origStyle
$ws.Range("E47").Value = '  +4.09%  '
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.40'
$c.Style = $origStyle
$ws.Range("E48").Value = '  +2.06%  '
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '146.26'
$c.Style = $origStyle
$ws.Range("E49").Value = '  +3.83%  '
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.15'
$c.Style = $origStyle
$ws.Range("E50").Value = '  -2.73%  '
$ws.Range("E51").Value = '  +2.20%  '
